# Apply the commit: insert two new weekly price records for "Perejil" right
# after the existing row 188 (i.e. as the new rows 189-190), pushing all the
# subsequent records down by two rows (old 189-322 become 191-324).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 189:190; Excel shifts rows 189:322 down to 191:324
# and extends the used range/dimension to A1:R324 automatically.
$ws.Rows("189:190").Insert()

# --- New row 189 ---
$ws.Range("A189").Value = 10
$ws.Range("B189").Value = "Vega Modelo de Temuco"
$ws.Range("C189").Value = "La Araucanía"
$ws.Range("D189").Value = 44741
$ws.Range("E189").Value = 9
$ws.Range("F189").Value = 100112044
$ws.Range("G189").Value = "Perejil"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 40
$ws.Range("K189").Value = 3300
$ws.Range("L189").Value = 3300
$ws.Range("M189").Value = 3300
$ws.Range("N189").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O189").Value = "Provincia de Cautín"
$ws.Range("P189").Value = 1100
$ws.Range("Q189").Value = 3
$ws.Range("R189").Value = "Hortaliza"

# --- New row 190 ---
$ws.Range("A190").Value = 10
$ws.Range("B190").Value = "Vega Modelo de Temuco"
$ws.Range("C190").Value = "La Araucanía"
$ws.Range("D190").Value = 44741
$ws.Range("E190").Value = 9
$ws.Range("F190").Value = 100112044
$ws.Range("G190").Value = "Perejil"
$ws.Range("H190").Value = "Sin especificar"
$ws.Range("I190").Value = "Primera"
$ws.Range("J190").Value = 50
$ws.Range("K190").Value = 3300
$ws.Range("L190").Value = 3300
$ws.Range("M190").Value = 3300
$ws.Range("N190").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O190").Value = "Región Metropolitana"
$ws.Range("P190").Value = 1100
$ws.Range("Q190").Value = 3
$ws.Range("R190").Value = "Hortaliza"
